$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 4 (ece score / kr-vs-kp)
$ws.Range("C4").Value = 2.253
$ws.Range("D4").Value = 1.807
$ws.Range("E4").Value = 1.04
$ws.Range("F4").Value = 1.653
$ws.Range("G4").Value = 1.083
$ws.Range("H4").Value = 1.878

# Row 5 (ece score / mushroom)
$ws.Range("C5").Value = 1.351
$ws.Range("D5").Value = 1.882
$ws.Range("E5").Value = 1.115
$ws.Range("F5").Value = 1.63
$ws.Range("G5").Value = 1.228
$ws.Range("H5").Value = 1.957

# Row 6 (brier score loss / kr-vs-kp)
$ws.Range("C6").Value = 0.773
$ws.Range("D6").Value = 0.412
$ws.Range("E6").Value = 0.394
$ws.Range("F6").Value = 0.633
$ws.Range("G6").Value = 0.431
$ws.Range("H6").Value = 0.384

# Row 7 (brier score loss / mushroom)
$ws.Range("C7").Value = 0.959
$ws.Range("D7").Value = 0.428
$ws.Range("E7").Value = 0.412
$ws.Range("F7").Value = 0.674
$ws.Range("G7").Value = 0.457
$ws.Range("H7").Value = 0.398

$wb.Save()
